$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new message row (row 7) describing the ChangeLobbyChatStatus message
$ws.Range("A7").Value = "ChangeLobbyChatStatus"
$ws.Range("B7").Value = 1005
$ws.Range("C7").Value = "Manda il nuovo stato (abilitata/disabilitata) della chat"

# Update the active selection to match the edited workbook
$ws.Range("C8").Select()
